$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "26.870.02"
$ws.Cells.Item(2, 5).Value = "  -1.41%  "

$ws.Cells.Item(3, 4).Value = "1.807.65"
$ws.Cells.Item(3, 5).Value = "  -1.08%  "

$ws.Cells.Item(4, 4).Value = "1.000"
$ws.Cells.Item(4, 5).Value = "  -0.47%  "

$ws.Cells.Item(5, 4).Value = "310.28"
$ws.Cells.Item(5, 5).Value = "  -1.09%  "

$ws.Cells.Item(6, 5).Value = "  -0.38%  "

$ws.Cells.Item(7, 5).Value = "  +4.78%  "

$ws.Cells.Item(8, 4).Value = "0.3674"
$ws.Cells.Item(8, 5).Value = "  -0.97%  "

$ws.Cells.Item(9, 4).Value = "0.07379"
$ws.Cells.Item(9, 5).Value = "  +1.67%  "

$ws.Cells.Item(10, 4).Value = "0.8562"
$ws.Cells.Item(10, 5).Value = "  -0.99%  "

$ws.Cells.Item(11, 5).Value = "  -1.95%  "

$ws.Cells.Item(12, 4).Value = "1.806.69"
$ws.Cells.Item(12, 5).Value = "  -1.41%  "

$ws.Cells.Item(13, 4).Value = "6.609"
$ws.Cells.Item(13, 5).Value = "  -1.91%  "

$ws.Cells.Item(14, 4).Value = "92.52"
$ws.Cells.Item(14, 5).Value = "  +3.37%  "

$ws.Cells.Item(15, 4).Value = "5.315"
$ws.Cells.Item(15, 5).Value = "  -0.19%  "

$ws.Cells.Item(16, 4).Value = "0.07080"
$ws.Cells.Item(16, 5).Value = "  -0.13%  "

$ws.Cells.Item(17, 4).Value = "1.002"
$ws.Cells.Item(17, 5).Value = "  -0.34%  "

$ws.Cells.Item(18, 4).Value = "0.000008744"
$ws.Cells.Item(18, 5).Value = "  -1.48%  "

$ws.Cells.Item(19, 4).Value = "1.0000"
$ws.Cells.Item(19, 5).Value = "  -0.46%  "

$ws.Cells.Item(20, 5).Value = "  -1.60%  "

$ws.Cells.Item(21, 4).Value = "26.884.78"
$ws.Cells.Item(21, 5).Value = "  -1.74%  "

$ws.Cells.Item(22, 4).Value = "5.162"
$ws.Cells.Item(22, 5).Value = "  +0.34%  "

$ws.Cells.Item(23, 4).Value = "10.86"
$ws.Cells.Item(23, 5).Value = "  -0.53%  "

$ws.Cells.Item(24, 5).Value = "  -0.03%  "

$ws.Cells.Item(25, 4).Value = "151.78"
$ws.Cells.Item(25, 5).Value = "  -0.53%  "

$ws.Cells.Item(26, 4).Value = "18.50"
$ws.Cells.Item(26, 5).Value = "  +0.39%  "

$ws.Cells.Item(27, 4).Value = "2.182"
$ws.Cells.Item(27, 5).Value = "  -1.16%  "

$ws.Cells.Item(28, 4).Value = "5.215"
$ws.Cells.Item(28, 5).Value = "  -0.67%  "

$ws.Cells.Item(29, 4).Value = "116.61"
$ws.Cells.Item(29, 5).Value = "  -0.06%  "

$ws.Cells.Item(30, 4).Value = "0.08828"
$ws.Cells.Item(30, 5).Value = "  -0.40%  "

$ws.Cells.Item(31, 4).Value = "0.7550"
$ws.Cells.Item(31, 5).Value = "  -0.69%  "

$ws.Cells.Item(32, 4).Value = "1.175"
$ws.Cells.Item(32, 5).Value = "  -2.08%  "

$ws.Cells.Item(33, 4).Value = "2.929"
$ws.Cells.Item(33, 5).Value = "  +4.50%  "

$ws.Cells.Item(34, 4).Value = "4.462"
$ws.Cells.Item(34, 5).Value = "  -0.07%  "

$ws.Cells.Item(35, 4).Value = "0.9998"
$ws.Cells.Item(35, 5).Value = "  -0.47%  "

$ws.Cells.Item(36, 4).Value = "1.089"
$ws.Cells.Item(36, 5).Value = "  -2.87%  "

$ws.Cells.Item(37, 4).Value = "0.01970"
$ws.Cells.Item(37, 5).Value = "  -0.41%  "

$ws.Cells.Item(38, 4).Value = "0.05196"
$ws.Cells.Item(38, 5).Value = "  -1.37%  "

$ws.Cells.Item(39, 4).Value = "0.5335"
$ws.Cells.Item(39, 5).Value = "  +5.24%  "

$ws.Cells.Item(40, 4).Value = "2.867"
$ws.Cells.Item(40, 5).Value = "  +0.01%  "

$ws.Cells.Item(41, 4).Value = "7.012"
$ws.Cells.Item(41, 5).Value = "  -4.51%  "

$ws.Cells.Item(42, 5).Value = "  -0.52%  "

$ws.Cells.Item(43, 4).Value = "0.5184"
$ws.Cells.Item(43, 5).Value = "  +8.91%  "

$ws.Cells.Item(44, 4).Value = "8.442"
$ws.Cells.Item(44, 5).Value = "  -3.33%  "

$ws.Cells.Item(45, 4).Value = "1.980"
$ws.Cells.Item(45, 5).Value = "  +6.03%  "

$ws.Cells.Item(46, 4).Value = "10.53"
$ws.Cells.Item(46, 5).Value = "  -0.92%  "

$ws.Cells.Item(47, 4).Value = "105.48"
$ws.Cells.Item(47, 5).Value = "  -1.99%  "

$ws.Cells.Item(48, 4).Value = "1.682"
$ws.Cells.Item(48, 5).Value = "  +0.32%  "

$ws.Cells.Item(49, 4).Value = "0.9994"
$ws.Cells.Item(49, 5).Value = "  -0.48%  "

$ws.Cells.Item(50, 4).Value = "0.06339"
$ws.Cells.Item(50, 5).Value = "  -0.89%  "

$ws.Cells.Item(51, 4).Value = "0.9227"
$ws.Cells.Item(51, 5).Value = "  +0.71%  "
